$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F ("想去人数" / wish-to-go count)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 594
$wsExpo.Range("F3").Value = 200
$wsExpo.Range("F4").Value = 440
$wsExpo.Range("F5").Value = 480
$wsExpo.Range("F6").Value = 278
$wsExpo.Range("F7").Value = 2520
$wsExpo.Range("F8").Value = 431
$wsExpo.Range("F9").Value = 6739
$wsExpo.Range("F11").Value = 431
$wsExpo.Range("F12").Value = 31

# Sheet "全部类型" (All types) - same underlying data merged with "演出"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 594
$wsAll.Range("F3").Value = 200
$wsAll.Range("F4").Value = 440
$wsAll.Range("F5").Value = 480
$wsAll.Range("F6").Value = 278
$wsAll.Range("F9").Value = 2520
$wsAll.Range("F10").Value = 431
$wsAll.Range("F11").Value = 6739
$wsAll.Range("F13").Value = 431
$wsAll.Range("F16").Value = 31
